$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Style-switching cells: copy formats from a same-style template cell, then set content ---
$ws.Range("F14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 1
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 5
$ws.Range("C20").Value = 5
$ws.Range("C21").Value = 17
$ws.Range("C23").Value = 1
$ws.Range("C24").Value = 12
$ws.Range("C26").Value = 10
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("F14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 2
$ws.Range("F14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1
$ws.Range("D14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("D17").Value = 6
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("D21").Value = 19
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 7
$ws.Range("D25").Value = 2
$ws.Range("D26").Value = 8
$ws.Range("F14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("E17").Value = -66.666666666666
$ws.Range("E18").Value = 100
$ws.Range("E19").Value = -28.571428571428
$ws.Range("E20").Value = 0
$ws.Range("E21").Value = -10.526315789473
$ws.Range("E23").Value = 0
$ws.Range("E24").Value = 71.428571428571
$ws.Range("E25").Value = 0
$ws.Range("E26").Value = 25
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = 100
$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = 0
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 13
$ws.Range("F20").Value = 13
$ws.Range("F21").Value = 60
$ws.Range("F23").Value = 3
$ws.Range("F24").Value = 37
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 27
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = 2
$ws.Range("G16").Value = 10
$ws.Range("G17").Value = 19
$ws.Range("G18").Value = 5
$ws.Range("G19").Value = 27
$ws.Range("G20").Value = 13
$ws.Range("G21").Value = 75
$ws.Range("G23").Value = 8
$ws.Range("G24").Value = 43
$ws.Range("G25").Value = 9
$ws.Range("G26").Value = 25
$ws.Range("G28").Value = 1
$ws.Range("G29").Value = 3
$ws.Range("G30").Value = 2
$ws.Range("H16").Value = -20
$ws.Range("H17").Value = -10.526315789473
$ws.Range("H18").Value = 40
$ws.Range("H19").Value = -51.851851851851
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = -20
$ws.Range("H23").Value = -62.5
$ws.Range("H24").Value = -13.953488372093
$ws.Range("H25").Value = -55.555555555555
$ws.Range("H26").Value = 8
$ws.Range("H28").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I14").Value = 2
$ws.Range("I16").Value = 48
$ws.Range("I17").Value = 105
$ws.Range("I18").Value = 28
$ws.Range("I19").Value = 126
$ws.Range("I20").Value = 52
$ws.Range("I21").Value = 368
$ws.Range("I23").Value = 52
$ws.Range("I24").Value = 249
$ws.Range("I25").Value = 56
$ws.Range("I26").Value = 157
$ws.Range("I29").Value = 7
$ws.Range("I30").Value = 6
$ws.Range("J14").Value = 2
$ws.Range("J16").Value = 42
$ws.Range("J17").Value = 110
$ws.Range("J18").Value = 25
$ws.Range("J19").Value = 160
$ws.Range("J20").Value = 71
$ws.Range("J21").Value = 416
$ws.Range("J23").Value = 59
$ws.Range("J24").Value = 265
$ws.Range("J25").Value = 54
$ws.Range("J26").Value = 161
$ws.Range("J29").Value = 9
$ws.Range("J30").Value = 8
$ws.Range("K14").Value = 0
$ws.Range("K16").Value = 14.285714285714
$ws.Range("K17").Value = -4.545454545454
$ws.Range("K18").Value = 12
$ws.Range("K19").Value = -21.25
$ws.Range("K20").Value = -26.760563380281
$ws.Range("K21").Value = -11.538461538461
$ws.Range("K23").Value = -11.864406779661
$ws.Range("K24").Value = -6.037735849056
$ws.Range("K25").Value = 3.703703703703
$ws.Range("K26").Value = -2.484472049689
$ws.Range("K29").Value = -22.222222222222
$ws.Range("K30").Value = -25
$ws.Range("L14").Value = -50
$ws.Range("L15").Value = 16.666666666666
$ws.Range("L16").Value = -28.358208955223
$ws.Range("L17").Value = -1.869158878504
$ws.Range("L18").Value = 0
$ws.Range("L19").Value = 10.526315789473
$ws.Range("L20").Value = -35
$ws.Range("L21").Value = -9.35960591133
$ws.Range("L22").Value = -83.333333333333
$ws.Range("L23").Value = -16.129032258064
$ws.Range("L24").Value = -16.442953020134
$ws.Range("L25").Value = -29.113924050632
$ws.Range("L26").Value = -13.259668508287
$ws.Range("L27").Value = 12.5
$ws.Range("L28").Value = 33.333333333333
$ws.Range("L29").Value = -12.5
$ws.Range("L30").Value = -25
$ws.Range("M14").Value = -33.333333333333
$ws.Range("M16").Value = -46.067415730337
$ws.Range("M17").Value = 41.891891891891
$ws.Range("M18").Value = -78.461538461538
$ws.Range("M19").Value = 17.757009345794
$ws.Range("M20").Value = -1.88679245283
$ws.Range("M21").Value = -20.518358531317
$ws.Range("M22").Value = -80
$ws.Range("M23").Value = 136.363636363636
$ws.Range("M24").Value = -2.734375
$ws.Range("M26").Value = -34.033613445378
$ws.Range("M29").Value = -36.363636363636
$ws.Range("M30").Value = -25
$ws.Range("N14").Value = -50
$ws.Range("N15").Value = -41.666666666666
$ws.Range("N16").Value = -84.364820846905
$ws.Range("N17").Value = -10.25641025641
$ws.Range("N18").Value = -91.054313099041
$ws.Range("N19").Value = -33.333333333333
$ws.Range("N20").Value = -95.042897998093
$ws.Range("N21").Value = -81.51682571572
$ws.Range("N29").Value = -50
$ws.Range("N30").Value = -53.846153846153

$excel.CutCopyMode = 0

